$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.597.17"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "2.295.89"
$ws.Range("E3").Value = "  +0.76%  "

$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.14"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  -1.37%  "

$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.43"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  -0.95%  "

$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.506"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +0.09%  "

$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -1.93%  "

$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.35"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  -3.09%  "

$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.98"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  +4.06%  "

$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0786"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  -0.67%  "

$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.118"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  +0.12%  "

$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "2.663.98"
$ws.Range("E15").Value = "  +1.11%  "

$ws.Range("D16").Value = "2.285.73"
$ws.Range("E16").Value = "  +0.69%  "

$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "42.561.42"
$ws.Range("E18").Value = "  -0.04%  "

$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.11"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  -6.77%  "

$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -0.93%  "

$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.98"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  -0.47%  "

$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.64"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +0.75%  "

$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +6.56%  "

$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.83"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("E25").Value = "  +0.07%  "

$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  -1.90%  "

$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.25"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  -3.39%  "

$ws.Range("E28").Value = "  +14.61%  "

$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.64"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = "  -0.30%  "

$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.02"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  -0.23%  "

$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.73"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  -4.04%  "

$ws.Range("E32").Value = "  +0.06%  "

$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.96"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  +0.01%  "

$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.47"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  -1.19%  "

$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.41"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  -7.26%  "

$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0696"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +1.00%  "

$ws.Range("E37").Value = "  -2.58%  "

$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0996"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  -1.12%  "

$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.74"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  +0.04%  "

$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.108"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  -1.34%  "

$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  -0.53%  "

$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.96"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  +10.34%  "

$ws.Range("D43").Value = "1.961.28"
$ws.Range("E43").Value = "  -2.04%  "

$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.46"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +4.76%  "

$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0278"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  -0.40%  "

$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  -0.66%  "

$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").Value = "2.526.32"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("E49").Value = "  -1.32%  "

$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.96"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  -1.09%  "

$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.20"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  -0.05%  "

